$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) cells we touch stay text, matching source inlineStr typing,
# since Excel would otherwise auto-convert numeric-looking strings to numbers.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.845.37'
$ws.Range('E2').Value = '  +1.22%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.838.95'
$ws.Range('E3').Value = '  +1.46%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '308.90'
$ws.Range('E5').Value = '  +1.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.004'
$ws.Range('E6').Value = '  +0.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4695'
$ws.Range('E7').Value = '  +3.33%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3655'
$ws.Range('E8').Value = '  +1.67%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07139'
$ws.Range('E9').Value = '  +0.51%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9207'
$ws.Range('E10').Value = '  +3.12%  '
$ws.Range('B11').Value = 'WrappedEther'
$ws.Range('C11').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '2.000.05'
$ws.Range('E11').Value = '  +10.96%  '
$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '19.53'
$ws.Range('E12').Value = '  +1.28%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07648'
$ws.Range('E13').Value = '  -1.00%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.281'
$ws.Range('E14').Value = '  +0.50%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.388'
$ws.Range('E15').Value = '  +1.39%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '88.06'
$ws.Range('E16').Value = '  +2.49%  '
$ws.Range('E17').Value = '  +0.18%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008628'
$ws.Range('E18').Value = '  +0.95%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.004'
$ws.Range('E19').Value = '  +0.13%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '26.891.05'
$ws.Range('E20').Value = '  +1.22%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.45'
$ws.Range('E21').Value = '  +2.25%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.003'
$ws.Range('E22').Value = '  +1.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.59'
$ws.Range('E23').Value = '  +0.99%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.922'
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.46'
$ws.Range('E25').Value = '  -0.39%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '18.19'
$ws.Range('E26').Value = '  +2.29%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.004'
$ws.Range('E27').Value = '  -1.00%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '114.00'
$ws.Range('E28').Value = '  +1.63%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.874'
$ws.Range('E29').Value = '  +1.13%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08815'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.207'
$ws.Range('E31').Value = '  +2.36%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.176'
$ws.Range('E32').Value = '  +6.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7459'
$ws.Range('E33').Value = '  +0.81%  '
$ws.Range('B34').Value = 'RenderToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.751'
$ws.Range('E34').Value = '  +1.21%  '
$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.472'
$ws.Range('E35').Value = '  +1.11%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.085'
$ws.Range('E36').Value = '  +1.30%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01939'
$ws.Range('E37').Value = '  +0.29%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05206'
$ws.Range('E38').Value = '  +2.50%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.955'
$ws.Range('E39').Value = '  +1.66%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5196'
$ws.Range('E40').Value = '  +2.26%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.945'
$ws.Range('E41').Value = '  +2.22%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1510'
$ws.Range('E42').Value = '  +0.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.144'
$ws.Range('E43').Value = '  +1.55%  '
$ws.Range('E44').Value = '  +5.36%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4690'
$ws.Range('E45').Value = '  +0.25%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.005'
$ws.Range('E46').Value = '  +0.17%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '101.33'
$ws.Range('E47').Value = '  +2.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.593'
$ws.Range('E48').Value = '  +1.93%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '64.83'
$ws.Range('E49').Value = '  +1.65%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06030'
$ws.Range('E50').Value = '  +0.64%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.8860'
$ws.Range('E51').Value = '  +4.81%  '
